$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> @{C=...; D=...; E=...} (only keys that change are set)
$changes = @{
    41  = @{ C = 126945; E = 662707757 }
    48  = @{ C = 150636; E = 275744616 }
    67  = @{ C = 27105;  E = 168714702 }
    69  = @{ C = 17891;  E = 103836984 }
    72  = @{ C = 331325; E = 635436311 }
    81  = @{ C = 88356;  E = 499673318 }
    104 = @{ C = 319351; D = 67905;  E = 561273145 }
    121 = @{ C = 1306327; D = 220386; E = 2275257402 }
    129 = @{ C = 633686; E = 3433172989 }
    132 = @{ C = 585952; E = 3470848013 }
    186 = @{ C = 236835; E = 1189990205 }
    189 = @{ C = 100471; E = 556436057 }
    204 = @{ C = 265654; E = 1271493760 }
    237 = @{ C = 283322; E = 1438427690 }
    240 = @{ C = 205926; E = 1069823635 }
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
